# Swap the contents of columns E (category-name) and F (group-code)
# for every row in the sheet (header row included), matching the
# commit's fix that reordered these two columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

$rngE = $ws.Range("E1:E$lastRow")
$rngF = $ws.Range("F1:F$lastRow")

$valuesE = $rngE.Value2
$valuesF = $rngF.Value2

$rngE.Value2 = $valuesF
$rngF.Value2 = $valuesE
